# Auto-generated Excel COM-interop script to update the cryptos worksheet.
# Applies the scraped price/volume refresh: updates the Price (column D) and
# Volume(1h) (column E) values for every data row (2-51), and also swaps the
# TRON / WrappedEther rows (11 and 12) which changed coin/link/price/volume.
#
# Column D values are numeric-looking strings (e.g. '0.9986', '30.643.01')
# that must stay TEXT (as in the source data) rather than become real numbers.
# We force text by prefixing with an apostrophe (Excel's own 'store as text'
# convention) and then reset the cell style back to Normal so no unintended
# number-format style gets attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'30.643.01"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +1.61%  '
$ws.Cells.Item(3, 4).Value = "'1.864.10"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +1.54%  '
$ws.Cells.Item(4, 4).Value = "'0.9986"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.19%  '
$ws.Cells.Item(5, 4).Value = "'236.26"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.22%  '
$ws.Cells.Item(6, 4).Value = "'0.9985"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.18%  '
$ws.Cells.Item(7, 4).Value = "'0.4743"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +1.60%  '
$ws.Cells.Item(8, 4).Value = "'0.2810"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +4.00%  '
$ws.Cells.Item(9, 4).Value = "'0.06455"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +2.98%  '
$ws.Cells.Item(10, 4).Value = "'18.44"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +15.13%  '
$ws.Cells.Item(11, 2).Value = 'WrappedEther'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(11, 4).Value = "'1.864.82"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +1.58%  '
$ws.Cells.Item(12, 2).Value = 'TRON'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(12, 4).Value = "'0.07474"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +1.09%  '
$ws.Cells.Item(13, 4).Value = "'91.04"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +8.94%  '
$ws.Cells.Item(14, 4).Value = "'5.035"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +2.31%  '
$ws.Cells.Item(15, 4).Value = "'0.6451"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +4.54%  '
$ws.Cells.Item(16, 4).Value = "'297.25"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +30.27%  '
$ws.Cells.Item(17, 4).Value = "'30.607.75"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +1.75%  '
$ws.Cells.Item(18, 4).Value = "'0.9988"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -0.15%  '
$ws.Cells.Item(19, 4).Value = "'12.93"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +4.74%  '
$ws.Cells.Item(20, 4).Value = "'0.000007444"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +2.38%  '
$ws.Cells.Item(21, 4).Value = "'2.102.27"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +1.32%  '
$ws.Cells.Item(22, 4).Value = "'0.9970"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.24%  '
$ws.Cells.Item(23, 4).Value = "'5.162"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +6.39%  '
$ws.Cells.Item(24, 4).Value = "'6.061"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +3.96%  '
$ws.Cells.Item(25, 4).Value = "'168.38"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +2.10%  '
$ws.Cells.Item(26, 4).Value = "'9.144"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -0.43%  '
$ws.Cells.Item(27, 4).Value = "'19.51"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +10.27%  '
$ws.Cells.Item(28, 4).Value = "'1.945"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +3.91%  '
$ws.Cells.Item(29, 4).Value = "'0.1041"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +1.02%  '
$ws.Cells.Item(30, 4).Value = "'1.345"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -1.79%  '
$ws.Cells.Item(31, 4).Value = "'4.090"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +0.44%  '
$ws.Cells.Item(32, 4).Value = "'3.930"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +3.82%  '
$ws.Cells.Item(33, 4).Value = "'0.04920"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +2.89%  '
$ws.Cells.Item(34, 4).Value = "'1.167"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +2.98%  '
$ws.Cells.Item(35, 4).Value = "'0.7155"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +0.91%  '
$ws.Cells.Item(36, 4).Value = "'2.704"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -0.10%  '
$ws.Cells.Item(37, 4).Value = "'0.01916"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +2.62%  '
$ws.Cells.Item(38, 4).Value = "'2.706"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +2.26%  '
$ws.Cells.Item(39, 4).Value = "'2.033"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +5.33%  '
$ws.Cells.Item(40, 4).Value = "'0.8871"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -0.67%  '
$ws.Cells.Item(41, 4).Value = "'106.83"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +2.51%  '
$ws.Cells.Item(42, 4).Value = "'0.9981"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -0.29%  '
$ws.Cells.Item(43, 4).Value = "'0.4168"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +4.21%  '
$ws.Cells.Item(44, 4).Value = "'5.531"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +0.22%  '
$ws.Cells.Item(45, 4).Value = "'7.313"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +5.35%  '
$ws.Cells.Item(46, 4).Value = "'64.14"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +7.50%  '
$ws.Cells.Item(47, 4).Value = "'0.1218"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +2.54%  '
$ws.Cells.Item(48, 4).Value = "'34.61"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +6.52%  '
$ws.Cells.Item(49, 4).Value = "'8.724"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +2.10%  '
$ws.Cells.Item(50, 4).Value = "'1.383"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +1.93%  '
$ws.Cells.Item(51, 4).Value = "'0.05526"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +0.41%  '
